$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17. This shifts the existing rows 17-68
# down to rows 18-69 (and carries the D-column date style/format along
# with them), matching the diff's row renumbering.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with this week's new data point
# (a fresh weekly price observation for Cilantro).
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44607
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112040
$ws.Range("G17").Value = "Cilantro"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 1800
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = 1900
$ws.Range("N17").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 950
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = "Hortaliza"
